# edit.ps1 -- applies the 0744 ordenanza reformatting described by the diff.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Paragraph 1: "Yerba Buena, 29 de Abril de 1996"
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(1)
$p.Format.KeepWithNext = 1
$p.Format.SpaceAfter = 12

# ---------------------------------------------------------------------
# Paragraph 2: "ORDENANZA Nº 744"
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(2)
$p.Format.KeepWithNext = 1
$p.Format.SpaceBefore = 12
$p.Format.SpaceAfter = 18
$p.Range.Bold = 1

# ---------------------------------------------------------------------
# Paragraph 3: "VISTO: La reciente implementación ..."
# Split into its own "VISTO: " paragraph plus a following paragraph.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("VISTO: ", $false, $false, $false, $false, $false, $true, 1, $false, "VISTO: `r", 2) | Out-Null

$pVisto = $d.Paragraphs.Item(3)
$pVisto.Alignment = 0
$pVisto.Format.KeepWithNext = 1
$pVisto.Format.SpaceBefore = 12
$pVisto.Format.SpaceAfter = 6
$pVisto.Range.Bold = 1

$pVistoBody = $d.Paragraphs.Item(4)
$pVistoBody.Alignment = 0
$pVistoBody.Format.KeepWithNext = 1
$pVistoBody.Format.SpaceAfter = 6
$insPt = $pVistoBody.Range
$insPt.Collapse(1)
$insPt.InsertAfter(" ")

# ---------------------------------------------------------------------
# Paragraph 5 (now): "CONSIDERANDO: Que los mismos fueron ..."
# Split into its own "CONSIDERANDO: " paragraph plus a following one.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("CONSIDERANDO: ", $false, $false, $false, $false, $false, $true, 1, $false, "CONSIDERANDO: `r", 2) | Out-Null

$pConsid = $d.Paragraphs.Item(5)
$pConsid.Alignment = 0
$pConsid.Format.KeepWithNext = 1
$pConsid.Format.SpaceBefore = 12
$pConsid.Format.SpaceAfter = 6
$pConsid.Range.Bold = 1

$pConsidBody = $d.Paragraphs.Item(6)
$pConsidBody.Alignment = 0
$pConsidBody.Format.KeepWithNext = 1
$pConsidBody.Format.SpaceAfter = 6
$insPt2 = $pConsidBody.Range
$insPt2.Collapse(1)
$insPt2.InsertAfter(" ")

# ---------------------------------------------------------------------
# Paragraphs 7-9 (now): "Que para el traspaso...", "Que la situación...",
# "Que se hace necesario..." -- only paragraph formatting changes.
# ---------------------------------------------------------------------
foreach ($idx in 7,8,9) {
    $p = $d.Paragraphs.Item($idx)
    $p.Alignment = 0
    $p.Format.KeepWithNext = 1
    $p.Format.SpaceAfter = 6
}

# ---------------------------------------------------------------------
# Paragraph 10 (now): "POR EL CONCEJO DELIBERANTE SANCIONA ..."
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("POR ", $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

$p = $d.Paragraphs.Item(10)
$p.Format.KeepWithNext = 1
$p.Format.SpaceBefore = 18
$p.Format.SpaceAfter = 18
$p.Format.LeftIndent = 99.2
$p.Format.RightIndent = 99.2
$p.Range.Bold = 1

# the trailing "." must stay non-bold
$rng = $d.Content
$rng.Find.Execute("EL CONCEJO DELIBERANTE SANCIONA CON FUERZA DE ORDENANZA.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$periodRng = $d.Range($rng.End - 1, $rng.End)
$periodRng.Bold = 0

# ---------------------------------------------------------------------
# Paragraph 11 (now): "ARTICULO PRIMERO: INCORPORASE ..."
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(11)
$p.Alignment = 0
$p.Format.KeepWithNext = 1
$p.Format.SpaceAfter = 6

$rng = $d.Content
$rng.Find.Execute("ARTICULO PRIMERO:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Font.Underline = 1

$spaces20 = "                    "
$rng = $d.Content
$rng.Find.Execute($spaces20 + "(", $false, $false, $false, $false, $false, $true, 1, $false, " (", 2) | Out-Null
$rng = $d.Content
$rng.Find.Execute($spaces20 + "(", $false, $false, $false, $false, $false, $true, 1, $false, " (", 2) | Out-Null

# ---------------------------------------------------------------------
# Paragraph 12 (now): "ARTICULO SEGUNDO: El departamento ..."
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(12)
$p.Alignment = 0
$p.Format.KeepWithNext = 1
$p.Format.SpaceAfter = 6

$rng = $d.Content
$rng.Find.Execute("ARTICULO SEGUNDO:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Font.Underline = 1

# ---------------------------------------------------------------------
# Paragraph 13 (now): "ARTICULO TERCERO: COMUNIQUESE, COPIESE Y ARCHIVESE."
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(13)
$p.Alignment = 0
$p.Format.KeepWithNext = 1
$p.Format.SpaceAfter = 6

$rng = $d.Content
$rng.Find.Execute("ARTICULO TERCERO:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Font.Underline = 1

